$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 45

# Copy the formatting (styles) from the last existing data row so the new
# row matches the workbook's established per-column styling (bold/border
# style on column A, date-time number format on column E) without minting
# brand-new style entries.
$ws.Range("A44:V44").Copy()
$ws.Range("A45:V45").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = 44
$ws.Cells.Item($row, 2).Value = "moldova"
$ws.Cells.Item($row, 3).Value = "super-liga"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45235.66666666666
$ws.Cells.Item($row, 6).Value = "Sheriff Tiraspol"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Milsami"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.22
$ws.Cells.Item($row, 11).Value = "04/11/2023 04:12"
$ws.Cells.Item($row, 12).Value = 1.25
$ws.Cells.Item($row, 13).Value = "05/11/2023 14:17"
$ws.Cells.Item($row, 14).Value = 4.88
$ws.Cells.Item($row, 15).Value = "04/11/2023 04:12"
$ws.Cells.Item($row, 16).Value = 5.02
$ws.Cells.Item($row, 17).Value = "05/11/2023 15:48"
$ws.Cells.Item($row, 18).Value = 7.92
$ws.Cells.Item($row, 19).Value = "04/11/2023 04:12"
$ws.Cells.Item($row, 20).Value = 9.619999999999999
$ws.Cells.Item($row, 21).Value = "05/11/2023 15:48"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/moldova/super-liga/sheriff-tiraspol-milsami/YVs8OrJe/"
